# "Added feature to copy lists and move to different board"
#
# The Week-N sheets in this Trello-board-style workbook each start with a
# header row ("Week 1 Requirements", "Week 2 Requirements", "Week 3
# Requirements") that duplicated the board/sheet name. That header row is
# removed from the first three week sheets (the rows below shift up), and
# the now-unused header strings disappear from the shared string table.

$wb = $excel.ActiveWorkbook

# Week 1: drop the "Week 1 Requirements" header row, everything shifts up.
$wsWeek1 = $wb.Worksheets.Item("Week 1")
$wsWeek1.Rows.Item(1).Delete()
$wsWeek1.Range("A1").EntireRow.Select()

# Week 2: drop the "Week 2 Requirements" header row, everything shifts up.
$wsWeek2 = $wb.Worksheets.Item("Week 2")
$wsWeek2.Rows.Item(1).Delete()
$wsWeek2.Range("A1").EntireRow.Select()

# Week 3: drop the "Week 3 Requirements" header row, everything shifts up.
$wsWeek3 = $wb.Worksheets.Item("Week 3")
$wsWeek3.Rows.Item(1).Delete()
$wsWeek3.Range("A1").EntireRow.Select()

# Week 4 is untouched by this change.

# Leave the workbook focused back on "Week Counter" (the tab that was active
# when the file was saved), with its last selected cell moved to D43.
$wsCounter = $wb.Worksheets.Item("Week Counter")
$wsCounter.Select()
$wsCounter.Range("D43").Select()
